$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$zhHandbackDate = "2016-03-31 05:22:34"
$deHandbackDate = "2016-03-31 05:22:49"

# --- Overview sheet: Status text refresh (same text used across all sheets) ---
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("B2").Value = $newStatus
$wsOv.Range("C2").Value = $newStatus
$wsOv.Range("B3").Value = $newStatus
$wsOv.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("H2").Value = $zhHandbackDate
$wsZh.Range("H3").Value = $zhHandbackDate

$wsZh.Range("F2").Value = "6339abdc-0900-4aad-8b02-6a01b73b3347.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/b382c276acc4038ddd245264c864a260ef8d5226/e2e/6339abdc-0900-4aad-8b02-6a01b73b3347.md", "", "", "6339abdc-0900-4aad-8b02-6a01b73b3347.md") | Out-Null

$wsZh.Range("G2").Value = "6339abdc-0900-4aad-8b02-6a01b73b3347.df65e2d6064421aff35691425e386963c2d0ed4f.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bb577ec4879f8a254cb623b2fab83eca983e2cb9/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/6339abdc-0900-4aad-8b02-6a01b73b3347.df65e2d6064421aff35691425e386963c2d0ed4f.zh-cn.xlf", "", "", "6339abdc-0900-4aad-8b02-6a01b73b3347.df65e2d6064421aff35691425e386963c2d0ed4f.zh-cn.xlf") | Out-Null

$wsZh.Range("F3").Value = "753516b8-1b0c-4f2a-a4fa-c621cde4422f.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/b382c276acc4038ddd245264c864a260ef8d5226/e2e/753516b8-1b0c-4f2a-a4fa-c621cde4422f.md", "", "", "753516b8-1b0c-4f2a-a4fa-c621cde4422f.md") | Out-Null

$wsZh.Range("G3").Value = "753516b8-1b0c-4f2a-a4fa-c621cde4422f.ede55c2c04de6316fff3ff9fc4f470d37a6b852f.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bb577ec4879f8a254cb623b2fab83eca983e2cb9/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/753516b8-1b0c-4f2a-a4fa-c621cde4422f.ede55c2c04de6316fff3ff9fc4f470d37a6b852f.zh-cn.xlf", "", "", "753516b8-1b0c-4f2a-a4fa-c621cde4422f.ede55c2c04de6316fff3ff9fc4f470d37a6b852f.zh-cn.xlf") | Out-Null

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("H2").Value = $deHandbackDate
$wsDe.Range("H3").Value = $deHandbackDate

$wsDe.Range("F2").Value = "6339abdc-0900-4aad-8b02-6a01b73b3347.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/b382c276acc4038ddd245264c864a260ef8d5226/e2e/6339abdc-0900-4aad-8b02-6a01b73b3347.md", "", "", "6339abdc-0900-4aad-8b02-6a01b73b3347.md") | Out-Null

$wsDe.Range("G2").Value = "6339abdc-0900-4aad-8b02-6a01b73b3347.df65e2d6064421aff35691425e386963c2d0ed4f.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfdca66c5a24df61f43d504561922fb2eca12093/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/6339abdc-0900-4aad-8b02-6a01b73b3347.df65e2d6064421aff35691425e386963c2d0ed4f.de-de.xlf", "", "", "6339abdc-0900-4aad-8b02-6a01b73b3347.df65e2d6064421aff35691425e386963c2d0ed4f.de-de.xlf") | Out-Null

$wsDe.Range("F3").Value = "753516b8-1b0c-4f2a-a4fa-c621cde4422f.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/b382c276acc4038ddd245264c864a260ef8d5226/e2e/753516b8-1b0c-4f2a-a4fa-c621cde4422f.md", "", "", "753516b8-1b0c-4f2a-a4fa-c621cde4422f.md") | Out-Null

$wsDe.Range("G3").Value = "753516b8-1b0c-4f2a-a4fa-c621cde4422f.ede55c2c04de6316fff3ff9fc4f470d37a6b852f.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfdca66c5a24df61f43d504561922fb2eca12093/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/753516b8-1b0c-4f2a-a4fa-c621cde4422f.ede55c2c04de6316fff3ff9fc4f470d37a6b852f.de-de.xlf", "", "", "753516b8-1b0c-4f2a-a4fa-c621cde4422f.ede55c2c04de6316fff3ff9fc4f470d37a6b852f.de-de.xlf") | Out-Null

Write-Host "Report generated for handback."
